$d = $word.ActiveDocument

# Update the date heading
$d.Content.Find.Execute("2024-12-25 Wednesday", $true, $false, $false, $false, $false, $true, 1, $false, "2024-12-26 Thursday", 2) | Out-Null

# Update the division-problem answer grid cell by cell (row, column) to
# avoid any ambiguity from duplicate strings appearing as both old and new values.
$t = $d.Tables.Item(1)

# Row 1
$t.Cell(1, 1).Range.Text = "82÷6=13, 4"
$t.Cell(1, 2).Range.Text = "29÷3=9, 2"
$t.Cell(1, 3).Range.Text = "72÷6=12, 0"
$t.Cell(1, 4).Range.Text = "22÷3=7, 1"
$t.Cell(1, 5).Range.Text = "55÷8=6, 7"

# Row 5
$t.Cell(5, 1).Range.Text = "76÷4=19, 0"
$t.Cell(5, 2).Range.Text = "62÷5=12, 2"
$t.Cell(5, 3).Range.Text = "93÷5=18, 3"
$t.Cell(5, 4).Range.Text = "95÷3=31, 2"
$t.Cell(5, 5).Range.Text = "36÷8=4, 4"

# Row 9
$t.Cell(9, 1).Range.Text = "42÷6=7, 0"
$t.Cell(9, 2).Range.Text = "38÷5=7, 3"
$t.Cell(9, 3).Range.Text = "22÷2=11, 0"
$t.Cell(9, 4).Range.Text = "95÷8=11, 7"
$t.Cell(9, 5).Range.Text = "18÷2=9, 0"

# Row 13
$t.Cell(13, 1).Range.Text = "27÷5=5, 2"
$t.Cell(13, 2).Range.Text = "46÷8=5, 6"
$t.Cell(13, 3).Range.Text = "78÷7=11, 1"
$t.Cell(13, 4).Range.Text = "80÷3=26, 2"
$t.Cell(13, 5).Range.Text = "19÷7=2, 5"

# Row 17
$t.Cell(17, 1).Range.Text = "82÷3=27, 1"
$t.Cell(17, 2).Range.Text = "12÷6=2, 0"
$t.Cell(17, 3).Range.Text = "40÷4=10, 0"
$t.Cell(17, 4).Range.Text = "57÷9=6, 3"
$t.Cell(17, 5).Range.Text = "29÷2=14, 1"
